$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 43; existing rows 43:113 shift down to 44:114
$ws.Rows.Item(43).Insert()

# Populate the new row 43 with the new record's data
$ws.Cells.Item(43, 1).Value = 9
$ws.Cells.Item(43, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(43, 3).Value = "Metropolitana"
$ws.Cells.Item(43, 4).Value = 44868
$ws.Cells.Item(43, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(43, 5).Value = 13
$ws.Cells.Item(43, 6).Value = 100114002
$ws.Cells.Item(43, 7).Value = "Camote"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 400
$ws.Cells.Item(43, 11).Value = 13000
$ws.Cells.Item(43, 12).Value = 13000
$ws.Cells.Item(43, 13).Value = 13000
$ws.Cells.Item(43, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(43, 15).Value = "Perú"
$ws.Cells.Item(43, 16).Value = 722
$ws.Cells.Item(43, 17).Value = 18
$ws.Cells.Item(43, 18).Value = "Hortaliza"
